$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A6").Value = "release/8.0.3"
$ws.Range("B6").Value = "X"
$ws.Range("C6").Value = "X"
$ws.Range("D6").Value = "X"
$ws.Range("E6").Value = "X"
